# Generate Report for Handback
# This script regenerates the localization-status report so that rows are
# re-sorted by (new) status, a new "Handed back: in sync with en-US" status
# appears for the two files that have been handed back, and the per-language
# detail sheets grow a "Latest Target File" / "Latest Handback File" column
# pair (F/G) that mirrors the existing source-file / handoff-file hyperlinks
# for the handed-back rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-21 08:17:58"

$ws1.Range("A3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-21 08:17:58"

$ws1.Range("A4").Value = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md"
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"
$ws1.Range("D4").Value = "2016-03-21 08:16:45"

$ws1.Range("A5").Value = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-21 08:17:58"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/614c999b-e180-4bc4-bd52-4849fba5c559.md", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/a49e2f31-4351-4b84-a105-cc480dc730a3.md", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md", "", "", "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d923ab19d85e6eb681bd8a4be5fdc15cb8d5004f/e2e/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md", "", "", "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-21 08:17:54"
$ws2.Range("F2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.md"
$ws2.Range("G2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-21 08:18:15"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-21 08:17:54"
$ws2.Range("F3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.md"
$ws2.Range("G3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-21 08:18:15"
$ws2.Range("J3").Value = "Include"

$ws2.Range("A4").Value = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "In Translation"
$ws2.Range("D4").Value = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-21 08:16:41"
$ws2.Range("F4").Value = ""
$ws2.Range("G4").Value = ""
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "Include"

$ws2.Range("A5").Value = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-21 08:17:54"
$ws2.Range("F5").Value = ""
$ws2.Range("G5").Value = ""
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("J5").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/614c999b-e180-4bc4-bd52-4849fba5c559.md", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce0f7565c3619d947004801fe0048bed7b382ad7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.zh-cn.xlf", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/614c999b-e180-4bc4-bd52-4849fba5c559.md", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce0f7565c3619d947004801fe0048bed7b382ad7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.zh-cn.xlf", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/a49e2f31-4351-4b84-a105-cc480dc730a3.md", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce0f7565c3619d947004801fe0048bed7b382ad7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.zh-cn.xlf", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/a49e2f31-4351-4b84-a105-cc480dc730a3.md", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce0f7565c3619d947004801fe0048bed7b382ad7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.zh-cn.xlf", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md", "", "", "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25e1da1c23593d68aa60eac33d03fe5ec62c421e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.zh-cn.xlf", "", "", "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d923ab19d85e6eb681bd8a4be5fdc15cb8d5004f/e2e/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md", "", "", "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce0f7565c3619d947004801fe0048bed7b382ad7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.zh-cn.xlf", "", "", "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-21 08:17:58"
$ws3.Range("F2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.md"
$ws3.Range("G2").Value = "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-21 08:18:21"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-21 08:17:58"
$ws3.Range("F3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.md"
$ws3.Range("G3").Value = "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-21 08:18:21"
$ws3.Range("J3").Value = "Include"

$ws3.Range("A4").Value = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "In Translation"
$ws3.Range("D4").Value = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-21 08:16:45"
$ws3.Range("F4").Value = ""
$ws3.Range("G4").Value = ""
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "Include"

$ws3.Range("A5").Value = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-21 08:17:58"
$ws3.Range("F5").Value = ""
$ws3.Range("G5").Value = ""
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("J5").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/614c999b-e180-4bc4-bd52-4849fba5c559.md", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d425be797fb598cac1f19a11aff57d68d27a66a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.de-de.xlf", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/614c999b-e180-4bc4-bd52-4849fba5c559.md", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d425be797fb598cac1f19a11aff57d68d27a66a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.de-de.xlf", "", "", "614c999b-e180-4bc4-bd52-4849fba5c559.b80cf98ff46335e4490e99cbb6bd716611246fc0.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/a49e2f31-4351-4b84-a105-cc480dc730a3.md", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d425be797fb598cac1f19a11aff57d68d27a66a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.de-de.xlf", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a704963cd5af42ecf3374ca6a81047f431433c9b/e2e/a49e2f31-4351-4b84-a105-cc480dc730a3.md", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d425be797fb598cac1f19a11aff57d68d27a66a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.de-de.xlf", "", "", "a49e2f31-4351-4b84-a105-cc480dc730a3.5926fec217692c2933f41351d5898731b78d200a.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md", "", "", "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daacd1925ee868f8937a97b32af7b49aedf83c28/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.de-de.xlf", "", "", "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d923ab19d85e6eb681bd8a4be5fdc15cb8d5004f/e2e/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md", "", "", "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d425be797fb598cac1f19a11aff57d68d27a66a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.de-de.xlf", "", "", "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.de-de.xlf")
